$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 94, pushing existing rows 94-124 down to 95-125.
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new record.
$ws.Range("A94").Value = 7
$ws.Range("B94").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C94").Value = 'Ñuble'
$ws.Range("D94").Value = 44468
$ws.Range("E94").Value = 16
$ws.Range("F94").Value = 100112006
$ws.Range("G94").Value = 'Repollo'
$ws.Range("H94").Value = 'Crespo record'
$ws.Range("I94").Value = 'Primera'
$ws.Range("J94").Value = 300
$ws.Range("K94").Value = 600
$ws.Range("L94").Value = 650
$ws.Range("M94").Value = 625
$ws.Range("N94").Value = '$/unidad'
$ws.Range("O94").Value = 'Provincia de Diguillín'
$ws.Range("P94").Value = 625
$ws.Range("Q94").Value = 1
$ws.Range("R94").Value = 'Hortaliza'
